$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: CO2 / fuel burn value (keeps its numeric-style xf, now holding text) ---
$ws.Range("A2").Value = "CO2"
$ws.Range("B2").Value = "3.15 kg (Kérosène)"

# --- Row 3: H2O ---
$ws.Range("A3").Value = "H2O"
$ws.Range("B3").Value = "1.237 kg "

# --- Row 4 & 5: SOx (renamed from SO2, now split into two rows, wrapped text) ---
$ws.Range("A4").Value = "SOx"
$ws.Range("B4").Value = "0.00084 kg - assuming 440 ppm"
$ws.Range("B4").WrapText = $true

$ws.Range("A5").Value = "SOx"
$ws.Range("B5").Value = "0.000114 kg – assuming 600 ppm"
$ws.Range("B5").WrapText = $true

# --- Row 6: NOx (new) ---
$ws.Range("A6").Value = "NOx"
$ws.Range("B6").Value = "0.0148 kg "

# --- Row 7: HC (new) ---
$ws.Range("A7").Value = "HC"
$ws.Range("B7").Value = "0.00032 kg"

# --- Row 8: CO (new) ---
$ws.Range("A8").Value = "CO"
$ws.Range("B8").Value = "0.00325 kg "

# --- Row 9: PM volatile (new) ---
$ws.Range("A9").Value = "PM volatile"
$ws.Range("B9").Value = "0.000092 kg "

# --- Column widths (closest achievable to the authored 10.66 / 31.33 "best fit" widths) ---
$ws.Columns.Item(1).ColumnWidth = 9.75
$ws.Columns.Item(2).ColumnWidth = 30.42

# --- Restore cursor/selection like the saved file ---
$ws.Range("F16").Select()
